# "added evaluation script results"
#
# The author logged a new entry in the time-tracking sheet:
#   2 hours spent "Schreiben von Evaluierungsscript und starten der Evaluierung"
# (writing the evaluation script and starting the evaluation), recorded in the
# first free row right above the "total" row. The SUM formula in column B
# then naturally picks up the extra 2 hours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B51").Value = 2
$ws.Range("C51").Value = "Schreiben von Evaluierungsscript und starten der Evaluierung"
